# "Generate Report for Handback"
#
# This localization-status report gets refreshed: every "Status" cell that
# used to read "In Translation" (the Overview sheet's per-language status
# columns, and the "Status" column on each language's own sheet) now reads
# "Handed back: in sync with en-US". The per-language "Latest Handback
# DateTime" stamps advance, the stale "handback file is not latest" Error
# Detail warning is cleared now that the handback is in sync, and a couple
# of report columns are resized so the (now longer) status / shorter error
# text reads better.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview: Status columns (zh-cn / de-de) now read "Handed back" ---
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack

# --- zh-cn: Status column, refreshed Latest Handback DateTime, clear stale Error Detail ---
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("L2").Value = "2017-02-15 05:57:40"
$zhcn.Range("R2").Value = ""

# --- de-de: Status column, refreshed Latest Handback DateTime, clear stale Error Detail ---
$dede.Range("C2").Value = $handedBack
$dede.Range("L2").Value = "2017-02-15 05:58:04"
$dede.Range("R2").Value = ""

# --- Column width tweaks (report layout refresh) ---
# Overview: widen the two status columns (E, F) to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: widen the Status column (C) and shrink Error Detail (R)
# now that it mostly holds short/blank text.
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(18).ColumnWidth = 12.833333333333332

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(18).ColumnWidth = 12.833333333333332
